$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "bAUep275"
$ws.Range("B2").Value = 23091337
$ws.Range("C2").Value = "dfqyfsf19"
$ws.Range("D2").Value = "wCbR92!%"
$ws.Range("F2").Value = "NOkBSCFB"
$ws.Range("G2").Value = "NhGG"

# Row 3
$ws.Range("A3").Value = "lIBwx541"
$ws.Range("B3").Value = 23091336
$ws.Range("C3").Value = "kpvdwjl80"
$ws.Range("D3").Value = "j!87#AFm"
$ws.Range("F3").Value = "gQJAUwQE"
$ws.Range("G3").Value = "acwp"

# Row 4
$ws.Range("A4").Value = "dVGiz463"
$ws.Range("B4").Value = 23091335
$ws.Range("C4").Value = "einfdpn36"
$ws.Range("D4").Value = "Z4ryN`$7%"
$ws.Range("F4").Value = "RAdaGpBc"
$ws.Range("G4").Value = "JQzP"
